$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data entries for "Betarraga" (rows 177-199) are being shifted down
# by one row (row N gets what row N-1 used to hold), a brand-new newest entry
# is written into row 177, and the block grows by one row: the old row 199
# becomes the new row 200.

# 1) Capture the existing block A177:R199 before any writes.
$srcRng = $ws.Range("A177:R199")
$vals = $srcRng.Value2

# 2) Shift that whole block down by one row -> A178:R200.
$dstRng = $ws.Range("A178:R200")
$dstRng.Value2 = $vals

# The newly-created row 200 needs the same date/time number format as the
# rest of column D (it was created fresh, so it has no inherited style).
$ws.Range("D200").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 3) Write the brand-new values for row 177 (the newest weekly entry).
$ws.Range("D177").Value = 44491
$ws.Range("J177").Value = 1200
$ws.Range("K177").Value = 1000
$ws.Range("L177").Value = 1200
$ws.Range("M177").Value = 1100
$ws.Range("P177").Value = 220
